$wb = $excel.ActiveWorkbook

# Remove Sheet2 entirely
$excel.DisplayAlerts = $false
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

# Work on Sheet1
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update the "App Type" column values from "dtc-model" to "model"
$ws1.Range("B2").Value = "model"
$ws1.Range("B3").Value = "model"
$ws1.Range("B4").Value = "model"

# Move the active selection from D22 to D14
$ws1.Range("D14").Select()
